# Final graphic tweaks and full testing
# Updates the lidargo config examples sheet:
#   - the "vol"/"vad" example column (C/D) now uses data level "a0" instead of
#     "a2" for its input file-name pattern, data_level_in, and the
#     corresponding min/max azimuth step values.
#   - leaves the final selection on D26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: file-name regex patterns per example column.
$ws.Range("C1").Value = "sa1.lidar.z05.a0.\d{8}.\d{6}.user5.nc"
$ws.Range("D1").Value = "sa1.lidar.z05.vad.a0.\d{8}.\d{6}.user5.nc"

# Row 7: data_level_in for the "vol" and "vad" example columns.
$ws.Range("C7").Value = "a0"
$ws.Range("D7").Value = "a0"

# Rows 10-11: min_azi_step / max_azi_step for the "vad" example column.
$ws.Range("D10").Value = -40
$ws.Range("D11").Value = 40

# Leave the sheet scrolled/selected the same way the author left it.
$ws.Range("D26").Select()
